$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on all touched cells so numeric-looking strings
# (e.g. "1.00", "69.354.45", "0.0000123") are preserved exactly as literal
# text instead of being re-interpreted/re-serialized as floating point numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.354.45'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.07%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.688.09'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.12%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '681.21'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.55%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '159.45'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.64%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.495'

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.01%  '

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -3.08%  '

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.40%  '

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -2.77%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.311.03'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.08%  '

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.50%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.679.80'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.41%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '69.323.81'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.17%  '

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.79%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '16.05'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.56%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.49'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.15%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '470.29'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -2.13%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.91'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.90%  '

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.13%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '80.01'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.21%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.835.66'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.07%  '

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.05%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000123'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -4.88%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.94'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -3.82%  '

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -3.40%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.69'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.92%  '

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -3.60%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.63'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.96%  '

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -3.12%  '

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.34%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.94'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.05%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.677.06'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.58%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.157'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -6.47%  '

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -2.14%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.28'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.75%  '

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -4.32%  '

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.07%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0907'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.03%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '170.14'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +4.19%  '

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.92%  '

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.83%  '

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '28.29'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -5.57%  '

$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'SuiNetwork'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.12'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.05%  '

$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'dogwifhat'
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.71'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.99%  '

$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'ONDO'
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.30'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -3.25%  '

$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'FLOKI'
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.000277'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.45%  '

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -3.25%  '
